$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-73 down to 25-74.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Cells.Item(24, 1).Value = 11
$ws.Cells.Item(24, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(24, 3).Value = "Bíobío"
$ws.Cells.Item(24, 4).Value = 44987
$ws.Cells.Item(24, 5).Value = 8
$ws.Cells.Item(24, 6).Value = 100112031
$ws.Cells.Item(24, 7).Value = "Poroto verde"
$ws.Cells.Item(24, 8).Value = "Magnum"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 220
$ws.Cells.Item(24, 11).Value = 12000
$ws.Cells.Item(24, 12).Value = 13000
$ws.Cells.Item(24, 13).Value = 12455
$ws.Cells.Item(24, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 16).Value = 498
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
